$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: "Linear Gesture" ---
$ws.Range("A14:K14").Copy()
$ws.Range("A15:K15").PasteSpecial(-4122)

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 1
$ws.Range("D15").Value = "A"
$ws.Range("E15").Value = "Linear Gesture"
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = "wait(3);`nvalidate1;`nlink_Click(EMML10_test_link);`nwait(2);`nvalidate2;`nlink_Click(gesture10_test_link);`nDrawGesture(linear_default);`nvalidate3;"
$ws.Range("H15").Value = "validate1`n{`nvalidate_PageTitle=Pocket Browser Tests`n};`nvalidate2`n{`nvalidate_PageTitle=EMML1.0 Test Index Page`n};`nvalidate3`n{`nvalidate_Result=Gesture detected for 1th time`nvalidate_Result=swipe`n};"

$ws.Rows.Item(15).RowHeight = 192

# --- Row 16: "Hold Gesture" ---
$ws.Range("A14:K14").Copy()
$ws.Range("A16:K16").PasteSpecial(-4122)

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 1
$ws.Range("D16").Value = "A"
$ws.Range("E16").Value = "Hold Gesture"
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = "wait(3);`nvalidate1;`nlink_Click(EMML10_test_link);`nwait(2);`nvalidate2;`nlink_Click(gesture10_test_link);`nDrawGesture(hold,100,100,6000);`nvalidate3;"
$ws.Range("H16").Value = "validate1`n{`nvalidate_PageTitle=Pocket Browser Tests`n};`nvalidate2`n{`nvalidate_PageTitle=EMML1.0 Test Index Page`n};`nvalidate3`n{`nvalidate_Result=Gesture detected for 3th time`nvalidate_Result=press`n};"

$ws.Rows.Item(16).RowHeight = 192

# --- View state: scroll so row 13 is at top, select A15 ---
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("A15").Select()
